$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 50
$ws.Range("F5").Value = 4600
$ws.Range("F6").Value = 1832
$ws.Range("F7").Value = 133
$ws.Range("F9").Value = 3090
$ws.Range("F12").Value = 254
$ws.Range("F13").Value = 617
$ws.Range("F14").Value = 528
$ws.Range("F15").Value = 520
$ws.Range("F16").Value = 370
$ws.Range("F17").Value = 131
$ws.Range("F20").Value = 120
$ws.Range("F21").Value = 1579
$ws.Range("F23").Value = 608
$ws.Range("F24").Value = 44
$ws.Range("F25").Value = 530
$ws.Range("F27").Value = 48
$ws.Range("F30").Value = 14
$ws.Range("F32").Value = 3711
$ws.Range("F33").Value = 755
$ws.Range("F35").Value = 542
$ws.Range("F37").Value = 1782

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 40

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 50
$ws.Range("F5").Value = 4600
$ws.Range("F6").Value = 1832
$ws.Range("F7").Value = 133
$ws.Range("F9").Value = 3090
$ws.Range("F12").Value = 254
$ws.Range("F13").Value = 617
$ws.Range("F14").Value = 528
$ws.Range("F15").Value = 520
$ws.Range("F17").Value = 370
$ws.Range("F18").Value = 131
$ws.Range("F21").Value = 120
$ws.Range("F22").Value = 1579
$ws.Range("F24").Value = 608
$ws.Range("F25").Value = 44
$ws.Range("F26").Value = 530
$ws.Range("F28").Value = 48
$ws.Range("F31").Value = 14
$ws.Range("F33").Value = 3711
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 755
$ws.Range("F37").Value = 542
$ws.Range("F39").Value = 1782
